$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 22: 25. Reverse Nodes in k-Group ----
$ws.Range("A22").Value = "25. Reverse Nodes in k-Group"

$ws.Range("B22").Value = "Hard"
$ws.Range("B22").Interior.Color = $ws.Range("B6").Interior.Color

$ws.Range("C22").Value = "Linked List"

$ws.Range("D22").Value = "We need a dummy node and to track kStart and kLast. Consider groups k at a time. If the kth node in the group is null, that is the break condition of the while loop. Use standard reverse code in a function and call iteratively."

$ws.Range("E22").Value = "https://leetcode.com/problems/reverse-nodes-in-k-group/solutions/11440/non-recursive-java-solution-and-idea/ "
$ws.Hyperlinks.Add($ws.Range("E22"), "https://leetcode.com/problems/reverse-nodes-in-k-group/solutions/11440/non-recursive-java-solution-and-idea/")
$ws.Range("E22").Style = "Hyperlink"

# ---- Row 23: 1448. Count Good Nodes in Binary Tree ----
$ws.Range("A23").Value = "1448. Count Good Nodes in Binary Tree"

$ws.Range("B23").Value = "Medium"
$ws.Range("B23").Interior.Color = $ws.Range("B2").Interior.Color

$ws.Range("C23").Value = "Trees"

$ws.Range("E23").Value = "https://leetcode.com/problems/count-good-nodes-in-binary-tree/solutions/635555/java-100-simple-easy-code-using-pre-order-tree-traversal/ "
$ws.Hyperlinks.Add($ws.Range("E23"), "https://leetcode.com/problems/count-good-nodes-in-binary-tree/solutions/635555/java-100-simple-easy-code-using-pre-order-tree-traversal/")
$ws.Range("E23").Style = "Hyperlink"

$ws.Range("D23").Value = "Straightforward, call DFS on root. Define DFS preorder function, but track max on the path."

# ---- Extend the table to include the new rows ----
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E23"))

# ---- Update the selected cell to reflect post-edit cursor position ----
[void]$ws.Range("D26").Select()
